$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update header cell B1 (still "Jira id" text-wise but original si node removed; value itself unchanged)
$ws.Range("B1").Value = "Jira id"

# Row 2
$ws.Range("C2").Value = "Verify that user is able to add an Article from ALL content search results page to a particular watchlist"

# Row 3
$ws.Range("C3").Value = "Verify that user is able to add a Patent from ALL content search results page to a particular watchlist"
$ws.Range("D3").Value = "N"
$ws.Range("E3").Value = "SKIP"

# Row 4
$ws.Range("C4").Value = "Verify that user is able to add a Post from ALL content search results page to a particular watchlist"
$ws.Range("D4").Value = "N"
$ws.Range("E4").Value = "SKIP"

# Clear rows 5-12 entirely
$ws.Range("A5:E12").Clear()
